# UploadBulkStudent.xlsx: add "password" and "TrackId" columns to the bulk
# student upload template so the importer can carry a login password and the
# student's track alongside the existing Name/Email/Mobile/... fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells, right after the existing "University" column (H).
$ws.Range("I1").Value = "password"
$ws.Range("J1").Value = "TrackId"

# Give the new "password" column (I) a custom width, matching the manual
# resize the author did for the header text. (Excel snaps ColumnWidth to
# whole-pixel increments on save, so this lands as close as possible to the
# template's shipped width.)
$ws.Columns.Item(9).ColumnWidth = 12.5

# The author's cursor ended up on the new password column when they saved.
$ws.Range("I8").Select()
